# Auto-generated Excel COM-interop script
# Applies numeric updates to the Coeurl_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H49").Value = 1177.1666
$ws.Range("I49").Value = 172.66667
$ws.Range("J49").Value = 2181.6667
$ws.Range("K49").Value = 518.00001
$ws.Range("L49").Value = 6545.000100000001
$ws.Range("M49").Value = -382.00001
$ws.Range("N49").Value = -6817.000100000001

$ws.Range("H51").Value = 2750
$ws.Range("J51").Value = 3200
$ws.Range("L51").Value = 3200
$ws.Range("N51").Value = -4168

$ws.Range("H96").Value = 410.06668
$ws.Range("I96").Value = 598.1111
$ws.Range("J96").Value = 128
$ws.Range("K96").Value = 1794.3333
$ws.Range("L96").Value = 384
$ws.Range("M96").Value = -421.3332999999998
$ws.Range("N96").Value = -3130

$ws.Range("H132").Value = 85027.375
$ws.Range("I132").Value = 47211.684
$ws.Range("K132").Value = 141635.052
$ws.Range("M132").Value = -139105.052

$ws.Range("H138").Value = 2929.16
$ws.Range("I138").Value = 1021.5
$ws.Range("J138").Value = 3671.0278
$ws.Range("K138").Value = 3064.5
$ws.Range("L138").Value = 11013.0834
$ws.Range("M138").Value = 2075.5
$ws.Range("N138").Value = -21293.0834


# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2288.4634
$ws.Range("I32").Value = 1888.4384
$ws.Range("K32").Value = 1888.4384
$ws.Range("M32").Value = -1601.4384

$ws.Range("H45").Value = 20167.215
$ws.Range("I45").Value = 17433.3
$ws.Range("J45").Value = 27002
$ws.Range("K45").Value = 17433.3
$ws.Range("L45").Value = 27002
$ws.Range("M45").Value = -17056.3
$ws.Range("N45").Value = -27756

$ws.Range("H122").Value = 2691
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550


# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H24").Value = 1250
$ws.Range("I24").Value = 1000
$ws.Range("J24").Value = 1500
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 1500
$ws.Range("M24").Value = -765
$ws.Range("N24").Value = -1970

$ws.Range("H94").Value = 1485.6511
$ws.Range("I94").Value = 1115.3422
$ws.Range("K94").Value = 1115.3422
$ws.Range("M94").Value = -664.3422


# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H45").Value = 45000
$ws.Range("I45").Value = 45000
$ws.Range("K45").Value = 45000
$ws.Range("M45").Value = -44407

$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50812

$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52808

$ws.Range("H134").Value = 25209.043
$ws.Range("I134").Value = 10779.368
$ws.Range("K134").Value = 32338.104
$ws.Range("M134").Value = -29803.104


# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H40").Value = 349.91666
$ws.Range("I40").Value = 211
$ws.Range("K40").Value = 844
$ws.Range("M40").Value = -775

$ws.Range("H86").Value = 428.83334
$ws.Range("I86").Value = 349.66666
$ws.Range("J86").Value = 508
$ws.Range("K86").Value = 1048.99998
$ws.Range("L86").Value = 1524
$ws.Range("M86").Value = 137.0000199999999
$ws.Range("N86").Value = -3896

$ws.Range("H89").Value = 428.83334
$ws.Range("I89").Value = 349.66666
$ws.Range("J89").Value = 508
$ws.Range("K89").Value = 3146.99994
$ws.Range("L89").Value = 4572
$ws.Range("M89").Value = 2781.00006
$ws.Range("N89").Value = -16428

$ws.Range("H97").Value = 1324.6
$ws.Range("I97").Value = 661.5
$ws.Range("J97").Value = 1766.6666
$ws.Range("K97").Value = 1984.5
$ws.Range("L97").Value = 5299.9998
$ws.Range("M97").Value = -1488.5
$ws.Range("N97").Value = -6291.9998

$ws.Range("H107").Value = 1056.5
$ws.Range("I107").Value = 921.1429000000001
$ws.Range("J107").Value = 2004
$ws.Range("K107").Value = 2763.4287
$ws.Range("L107").Value = 6012
$ws.Range("M107").Value = -843.4287000000004
$ws.Range("N107").Value = -9852

$ws.Range("H124").Value = 9997.799999999999
$ws.Range("J124").Value = 9997.799999999999
$ws.Range("L124").Value = 29993.4
$ws.Range("N124").Value = -39813.39999999999

$ws.Range("H129").Value = 667.4167
$ws.Range("I129").Value = 591.7273
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 1775.1819
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 3224.8181
$ws.Range("N129").Value = -14500

$ws.Range("H131").Value = 29277.422
$ws.Range("I131").Value = 333833
$ws.Range("J131").Value = 3172.6572
$ws.Range("K131").Value = 1001499
$ws.Range("L131").Value = 9517.971600000001
$ws.Range("M131").Value = -996459
$ws.Range("N131").Value = -19597.9716

$ws.Range("H132").Value = 1340.3125
$ws.Range("I132").Value = 1167
$ws.Range("J132").Value = 1721.6
$ws.Range("K132").Value = 10503
$ws.Range("L132").Value = 15494.4
$ws.Range("M132").Value = -7973
$ws.Range("N132").Value = -20554.4

$ws.Range("H137").Value = 2860.6924
$ws.Range("I137").Value = 2761.125
$ws.Range("J137").Value = 3020
$ws.Range("K137").Value = 8283.375
$ws.Range("L137").Value = 9060
$ws.Range("M137").Value = -3183.375
$ws.Range("N137").Value = -19260


# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H122").Value = 2906.5
$ws.Range("I122").Value = 2735.5454
$ws.Range("J122").Value = 3533.3333
$ws.Range("K122").Value = 8206.636200000001
$ws.Range("L122").Value = 10599.9999
$ws.Range("M122").Value = -5756.636200000001
$ws.Range("N122").Value = -15499.9999

$ws.Range("H132").Value = 479819.25
$ws.Range("I132").Value = 479819.25
$ws.Range("K132").Value = 1439457.75
$ws.Range("M132").Value = -1436927.75


# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 6250.5
$ws.Range("I40").Value = 5900.6
$ws.Range("K40").Value = 5900.6
$ws.Range("M40").Value = -5764.6

$ws.Range("H55").Value = 198.9375
$ws.Range("I55").Value = 195.57143
$ws.Range("J55").Value = 222.5
$ws.Range("K55").Value = 195.57143
$ws.Range("L55").Value = 222.5
$ws.Range("M55").Value = -22.57142999999999
$ws.Range("N55").Value = -568.5

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H76").Value = 38162
$ws.Range("I76").Value = 38000
$ws.Range("K76").Value = 38000
$ws.Range("M76").Value = -37662

$ws.Range("H79").Value = 38162
$ws.Range("I79").Value = 38000
$ws.Range("K79").Value = 38000
$ws.Range("M79").Value = -36830

$ws.Range("H93").Value = 465448.38
$ws.Range("I93").Value = 619659.5600000001
$ws.Range("J93").Value = 2814.8333
$ws.Range("K93").Value = 619659.5600000001
$ws.Range("L93").Value = 2814.8333
$ws.Range("M93").Value = -618411.5600000001
$ws.Range("N93").Value = -5310.8333

$ws.Range("H100").Value = 56546
$ws.Range("I100").Value = 67482.5
$ws.Range("J100").Value = 12800
$ws.Range("K100").Value = 67482.5
$ws.Range("L100").Value = 12800
$ws.Range("M100").Value = -66941.5
$ws.Range("N100").Value = -13882


# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H30").Value = 62508240
$ws.Range("J30").Value = 62508240
$ws.Range("L30").Value = 62508240
$ws.Range("N30").Value = -62508454

$ws.Range("H81").Value = 9355.9
$ws.Range("I81").Value = 27054.75
$ws.Range("J81").Value = 4931.1875
$ws.Range("K81").Value = 54109.5
$ws.Range("L81").Value = 9862.375
$ws.Range("M81").Value = -53048.5
$ws.Range("N81").Value = -11984.375

$ws.Range("H84").Value = 9355.9
$ws.Range("I84").Value = 27054.75
$ws.Range("J84").Value = 4931.1875
$ws.Range("K84").Value = 270547.5
$ws.Range("L84").Value = 49311.875
$ws.Range("M84").Value = -265243.5
$ws.Range("N84").Value = -59919.875

$ws.Range("H126").Value = 4831.5
$ws.Range("I126").Value = 4608.5
$ws.Range("K126").Value = 13825.5
$ws.Range("M126").Value = -11355.5

$ws.Range("H132").Value = 2456.6182
$ws.Range("I132").Value = 2710.8936
$ws.Range("J132").Value = 962.75
$ws.Range("K132").Value = 8132.6808
$ws.Range("L132").Value = 2888.25
$ws.Range("M132").Value = -5602.6808
$ws.Range("N132").Value = -7948.25

